$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 7 (Ano 2025) with refreshed billing data
$ws.Range("B7").Value = 2893276.92
$ws.Range("C7").Value = -34.88125029861055
$ws.Range("D7").Value = 2917
$ws.Range("E7").Value = 2917
$ws.Range("F7").Value = 991.8673020226259
$ws.Range("G7").Value = 5.725882271436555
